$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.421.71'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '2.545.97'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''312.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Value = '''100.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.61%  '
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '''0.528'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('D10').Value = '''36.17'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.33%  '
$ws.Range('D11').Value = '''0.0802'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').Value = '''7.38'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '2.940.44'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('E15').Value = '  +7.56%  '
$ws.Range('D16').Value = '2.581.09'
$ws.Range('E16').Value = '  +5.92%  '
$ws.Range('D17').Value = '''0.842'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '42.463.63'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').Value = '''6.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0951'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').Value = '''12.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('D22').Value = '''69.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').Value = '''243.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.45%  '
$ws.Range('D24').Value = '''2.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.38%  '
$ws.Range('D25').Value = '''2.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '''1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''26.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('D29').Value = '''40.19'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '''10.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('D31').Value = '''157.54'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.71%  '
$ws.Range('D32').Value = '''5.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.52%  '
$ws.Range('D33').Value = '''2.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.69%  '
$ws.Range('D34').Value = '''0.0800'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.30%  '
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('E36').Value = '  -3.24%  '
$ws.Range('E37').Value = '  -3.10%  '
$ws.Range('D38').Value = '''18.20'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.75%  '
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('D41').Value = '''4.18'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.63%  '
$ws.Range('D42').Value = '''21.95'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').Value = '''3.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.09%  '
$ws.Range('D45').Value = '''0.0298'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('D46').Value = '1.962.79'
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('D47').Value = '''8.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('D48').Value = '2.795.98'
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D49').Value = '''80.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.67%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '''0.192'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.16%  '
$ws.Range('D51').Value = '''72.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.02%  '
